$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix three cells that should read "NaN" (text) instead of a number ---
$ws.Range("CU17").Value = "NaN"
$ws.Range("CR31").Value = "NaN"
$ws.Range("AG173").Value = "NaN"

# --- Append new row 206 with the day's data (columns A:DX) ---
$row206 = New-Object "object[,]" 1,128
$row206[0,0] = 44100
$row206[0,1] = 806038
$row206[0,2] = 2742
$row206[0,3] = 110451
$row206[0,4] = 67216
$row206[0,5] = 262954
$row206[0,6] = 28523
$row206[0,7] = 6672
$row206[0,8] = 5515
$row206[0,9] = 8370
$row206[0,10] = 9098
$row206[0,11] = 18887
$row206[0,12] = 3989
$row206[0,13] = 23628
$row206[0,14] = 32438
$row206[0,15] = 7910
$row206[0,16] = 10754
$row206[0,17] = 15076
$row206[0,18] = 14563
$row206[0,19] = 18069
$row206[0,20] = 15323
$row206[0,21] = 3729
$row206[0,22] = 3267
$row206[0,23] = 10576
$row206[0,24] = 29903
$row206[0,25] = 13973
$row206[0,26] = 11881
$row206[0,27] = 60369
$row206[0,28] = 2167
$row206[0,29] = 1127
$row206[0,30] = 750
$row206[0,31] = 473
$row206[0,32] = 750
$row206[0,33] = 473
$row206[0,34] = 736
$row206[0,35] = 2060
$row206[0,36] = 5580
$row206[0,37] = 37962
$row206[0,38] = 9761
$row206[0,39] = 2560
$row206[0,40] = 46838
$row206[0,41] = 1106
$row206[0,42] = 22773
$row206[0,43] = 1526
$row206[0,44] = 10368
$row206[0,45] = 1669
$row206[0,46] = 1606
$row206[0,47] = 8265
$row206[0,48] = 2011
$row206[0,49] = 964
$row206[0,50] = 2501
$row206[0,51] = 2687
$row206[0,52] = 63932
$row206[0,53] = 14065
$row206[0,54] = 6445
$row206[0,55] = 9800
$row206[0,56] = 7103
$row206[0,57] = 257
$row206[0,58] = 1469
$row206[0,59] = 2731
$row206[0,60] = 744
$row206[0,61] = 2166
$row206[0,62] = 9830
$row206[0,63] = 9549
$row206[0,64] = 10601
$row206[0,65] = 14311
$row206[0,66] = 1967
$row206[0,67] = 904
$row206[0,68] = 13699
$row206[0,69] = 11217
$row206[0,70] = 13113
$row206[0,71] = 3001
$row206[0,72] = 2193
$row206[0,73] = 5872
$row206[0,74] = 4857
$row206[0,75] = 2325
$row206[0,76] = 5921
$row206[0,77] = 3715
$row206[0,78] = 2211
$row206[0,79] = 1013
$row206[0,80] = 3020
$row206[0,81] = 2243
$row206[0,82] = 2006
$row206[0,83] = 1804
$row206[0,84] = 6522
$row206[0,85] = 2185
$row206[0,86] = 1467
$row206[0,87] = 1801
$row206[0,88] = 2112
$row206[0,89] = 2210
$row206[0,90] = 2588
$row206[0,91] = 1736
$row206[0,92] = 1217
$row206[0,93] = 1223
$row206[0,94] = 998
$row206[0,95] = 3439
$row206[0,96] = 1478
$row206[0,97] = 967
$row206[0,98] = 1112
$row206[0,99] = 1736
$row206[0,100] = 1603
$row206[0,101] = 814
$row206[0,102] = 893
$row206[0,103] = 1335
$row206[0,104] = 1677
$row206[0,105] = 1552
$row206[0,106] = 1616
$row206[0,107] = 1250
$row206[0,108] = 334
$row206[0,109] = 371
$row206[0,110] = 832
$row206[0,111] = 774
$row206[0,112] = 495
$row206[0,113] = 544
$row206[0,114] = 384
$row206[0,115] = 670
$row206[0,116] = 757
$row206[0,117] = 527
$row206[0,118] = 494
$row206[0,119] = 374
$row206[0,120] = 524
$row206[0,121] = 137612
$row206[0,122] = 340983
$row206[0,123] = 19631
$row206[0,124] = 148042
$row206[0,125] = 91525
$row206[0,126] = 45275
$row206[0,127] = 12849
$ws.Range("A206:DX206").Value = $row206

# --- Move the selection to the new last cell, matching the author's saved view ---
$ws.Range("B206").Select()
